$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 203
$ws.Range("I2").Value = 203.66667
$ws.Range("K2").Value = 203.66667
$ws.Range("M2").Value = -90.66667000000001
$ws.Range("H28").Value = 4234.75
$ws.Range("I28").Value = 2466.6667
$ws.Range("J28").Value = 5295.6
$ws.Range("K28").Value = 2466.6667
$ws.Range("L28").Value = 5295.6
$ws.Range("M28").Value = -1981.6667
$ws.Range("N28").Value = -6265.6
$ws.Range("H33").Value = 749.3333
$ws.Range("I33").Value = 448.5
$ws.Range("K33").Value = 448.5
$ws.Range("M33").Value = -219.5
$ws.Range("H51").Value = 3000
$ws.Range("J51").Value = 3000
$ws.Range("L51").Value = 3000
$ws.Range("N51").Value = -3968
$ws.Range("H61").Value = 23350
$ws.Range("I61").Value = 23350
$ws.Range("K61").Value = 70050
$ws.Range("M61").Value = -69878
$ws.Range("H64").Value = 9750
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 9750
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 9750
$ws.Range("M64").ClearContents() | Out-Null
$ws.Range("N64").Value = -10246
$ws.Range("H67").Value = 9750
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 9750
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 9750
$ws.Range("M67").ClearContents() | Out-Null
$ws.Range("N67").Value = -11466
$ws.Range("H80").Value = 731.3333
$ws.Range("I80").Value = 740
$ws.Range("J80").Value = 714
$ws.Range("K80").Value = 2220
$ws.Range("L80").Value = 2142
$ws.Range("M80").Value = -1222
$ws.Range("N80").Value = -4138
$ws.Range("H83").Value = 731.3333
$ws.Range("I83").Value = 740
$ws.Range("J83").Value = 714
$ws.Range("K83").Value = 6660
$ws.Range("L83").Value = 6426
$ws.Range("M83").Value = -1668
$ws.Range("N83").Value = -16410
$ws.Range("H118").Value = 2733
$ws.Range("I118").Value = 1200
$ws.Range("K118").Value = 3600
$ws.Range("M118").Value = -1943
$ws.Range("H125").Value = 1865.6666
$ws.Range("I125").Value = 1823.75
$ws.Range("K125").Value = 16413.75
$ws.Range("M125").Value = -13953.75
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 931.7368
$ws.Range("I2").Value = 793.7143
$ws.Range("J2").Value = 1318.2
$ws.Range("K2").Value = 793.7143
$ws.Range("L2").Value = 1318.2
$ws.Range("M2").Value = -680.7143
$ws.Range("N2").Value = -1544.2
$ws.Range("H32").Value = 3034241.2
$ws.Range("I32").Value = 4096.048
$ws.Range("K32").Value = 4096.048
$ws.Range("M32").Value = -3809.048
$ws.Range("H61").Value = 3410.6365
$ws.Range("I61").Value = 2202.4666
$ws.Range("J61").Value = 5999.5713
$ws.Range("K61").Value = 2202.4666
$ws.Range("L61").Value = 5999.5713
$ws.Range("M61").Value = -1990.4666
$ws.Range("N61").Value = -6423.5713
$ws.Range("H116").Value = 931.7368
$ws.Range("I116").Value = 793.7143
$ws.Range("J116").Value = 1318.2
$ws.Range("K116").Value = 793.7143
$ws.Range("L116").Value = 1318.2
$ws.Range("M116").Value = 1500.2857
$ws.Range("N116").Value = -5906.2
$ws.Range("H132").Value = 3230.762
$ws.Range("I132").Value = 3187.7368
$ws.Range("K132").Value = 9563.2104
$ws.Range("M132").Value = -7033.2104
$ws.Range("H136").Value = 3410.6365
$ws.Range("I136").Value = 2202.4666
$ws.Range("J136").Value = 5999.5713
$ws.Range("K136").Value = 6607.399800000001
$ws.Range("L136").Value = 17998.7139
$ws.Range("M136").Value = -4057.399800000001
$ws.Range("N136").Value = -23098.7139
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 931.7368
$ws.Range("I3").Value = 793.7143
$ws.Range("J3").Value = 1318.2
$ws.Range("K3").Value = 793.7143
$ws.Range("L3").Value = 1318.2
$ws.Range("M3").Value = -679.7143
$ws.Range("N3").Value = -1546.2
$ws.Range("H99").Value = 1300.8
$ws.Range("I99").Value = 1187.125
$ws.Range("J99").Value = 1755.5
$ws.Range("K99").Value = 1187.125
$ws.Range("L99").Value = 1755.5
$ws.Range("M99").Value = 310.875
$ws.Range("N99").Value = -4751.5
$ws.Range("H134").Value = 1749.5
$ws.Range("I134").Value = 1749.5
$ws.Range("K134").Value = 5248.5
$ws.Range("M134").Value = -2713.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3850
$ws.Range("I16").Value = 1200
$ws.Range("J16").Value = 6500
$ws.Range("K16").Value = 1200
$ws.Range("L16").Value = 6500
$ws.Range("M16").Value = -913
$ws.Range("N16").Value = -7074
$ws.Range("H58").Value = 3616.1667
$ws.Range("I58").Value = 3155.6667
$ws.Range("K58").Value = 3155.6667
$ws.Range("M58").Value = -2952.6667
$ws.Range("H113").Value = 3850
$ws.Range("I113").Value = 1200
$ws.Range("J113").Value = 6500
$ws.Range("K113").Value = 1200
$ws.Range("L113").Value = 6500
$ws.Range("M113").Value = 970
$ws.Range("N113").Value = -10840
$ws.Range("H136").Value = 3616.1667
$ws.Range("I136").Value = 3155.6667
$ws.Range("K136").Value = 9467.000100000001
$ws.Range("M136").Value = -6917.000100000001
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents() | Out-Null
$ws.Range("H141").Value = 45000
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents() | Out-Null
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 226498.44
$ws.Range("I4").Value = 402198
$ws.Range("J4").Value = 6874
$ws.Range("K4").Value = 1206594
$ws.Range("L4").Value = 20622
$ws.Range("M4").Value = -1206482
$ws.Range("N4").Value = -20846
$ws.Range("H6").Value = 218.36363
$ws.Range("I6").Value = 44.666668
$ws.Range("K6").Value = 134.000004
$ws.Range("M6").Value = -21.00000399999999
$ws.Range("H12").Value = 147.21428
$ws.Range("J12").Value = 161.55556
$ws.Range("L12").Value = 484.66668
$ws.Range("N12").Value = -830.66668
$ws.Range("H14").Value = 538.25
$ws.Range("I14").Value = 538.25
$ws.Range("K14").Value = 1614.75
$ws.Range("M14").Value = -1441.75
$ws.Range("H98").Value = 501.66666
$ws.Range("I98").Value = 500.5
$ws.Range("J98").Value = 504
$ws.Range("K98").Value = 1501.5
$ws.Range("L98").Value = 1512
$ws.Range("M98").Value = -3.5
$ws.Range("N98").Value = -4508
$ws.Range("H107").Value = 507.64285
$ws.Range("I107").Value = 330.16666
$ws.Range("J107").Value = 640.75
$ws.Range("K107").Value = 990.4999799999999
$ws.Range("L107").Value = 1922.25
$ws.Range("M107").Value = 929.5000200000001
$ws.Range("N107").Value = -5762.25
$ws.Range("H113").Value = 527.9091
$ws.Range("I113").Value = 545.3333
$ws.Range("J113").Value = 449.5
$ws.Range("K113").Value = 1635.9999
$ws.Range("L113").Value = 1348.5
$ws.Range("M113").Value = 534.0001
$ws.Range("N113").Value = -5688.5
$ws.Range("H140").Value = 3127.25
$ws.Range("I140").Value = 2859.7144
$ws.Range("K140").Value = 8579.143199999999
$ws.Range("M140").Value = -3399.143199999999
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 7487360.5
$ws.Range("I11").Value = 6724073
$ws.Range("J11").Value = 11609113
$ws.Range("K11").Value = 6724073
$ws.Range("L11").Value = 11609113
$ws.Range("M11").Value = -6723934
$ws.Range("N11").Value = -11609391
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents() | Out-Null
$ws.Range("H102").Value = 1835.7778
$ws.Range("I102").Value = 1902.75
$ws.Range("K102").Value = 1902.75
$ws.Range("M102").Value = -280.75
$ws.Range("H132").Value = 2703.258
$ws.Range("I132").Value = 2368.0833
$ws.Range("J132").Value = 3852.4285
$ws.Range("K132").Value = 7104.249899999999
$ws.Range("L132").Value = 11557.2855
$ws.Range("M132").Value = -4574.249899999999
$ws.Range("N132").Value = -16617.2855
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H35").Value = 1906.2
$ws.Range("I35").Value = 1382.75
$ws.Range("K35").Value = 1382.75
$ws.Range("M35").Value = -1046.75
$ws.Range("H55").Value = 2254.0557
$ws.Range("I55").Value = 1134.875
$ws.Range("J55").Value = 3149.4
$ws.Range("K55").Value = 1134.875
$ws.Range("L55").Value = 3149.4
$ws.Range("M55").Value = -961.875
$ws.Range("N55").Value = -3495.4
$ws.Range("H132").Value = 1540.1428
$ws.Range("I132").Value = 1640.5
$ws.Range("K132").Value = 4921.5
$ws.Range("M132").Value = -2391.5
$ws.Range("H136").Value = 4625.25
$ws.Range("I136").Value = 3334
$ws.Range("J136").Value = 5400
$ws.Range("K136").Value = 10002
$ws.Range("L136").Value = 16200
$ws.Range("M136").Value = -7452
$ws.Range("N136").Value = -21300
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 8950
$ws.Range("J29").Value = 8950
$ws.Range("L29").Value = 8950
$ws.Range("N29").Value = -9530
$ws.Range("H96").Value = 1450
$ws.Range("I96").Value = 1483.3334
$ws.Range("J96").Value = 1250
$ws.Range("K96").Value = 1483.3334
$ws.Range("M96").Value = -110.3334
$ws.Range("N96").Value = -3996
$ws.Range("H132").Value = 2318.8572
$ws.Range("I132").Value = 2318.8572
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 6956.571599999999
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -4426.571599999999
$ws.Range("N132").ClearContents() | Out-Null
$ws.Range("H136").Value = 4157.273
$ws.Range("I136").Value = 3135
$ws.Range("J136").Value = 5179.5454
$ws.Range("K136").Value = 9405
$ws.Range("L136").Value = 15538.6362
$ws.Range("M136").Value = -6855
$ws.Range("N136").Value = -20638.6362

Write-Host "Updated 255 cells across 8 sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)"
